# Rename the inline logo pictures living in the document's headers and
# footers:
#   - BTec logo  (header story): image1.jpg -> image2.jpg
#   - Pearson logo (footer story): image2.png -> image1.png
#
# The document keeps each logo as an <wp:inline> picture placed directly
# in the primary + first-page header/footer stories, so walk every
# story range and touch the inline picture(s) we find there.

$d = $word.ActiveDocument

# Word header/footer "story type" constants:
#   6  = wdEvenPagesHeaderStory     7  = wdPrimaryHeaderStory
#   8  = wdEvenPagesFooterStory     9  = wdPrimaryFooterStory
#  10  = wdFirstPageHeaderStory    11  = wdFirstPageFooterStory
$headerStoryTypes = @(6, 7, 10)
$footerStoryTypes = @(8, 9, 11)

foreach ($story in $d.StoryRanges) {
    $storyType = $story.StoryType
    $isHeaderStory = $headerStoryTypes -contains $storyType
    $isFooterStory = $footerStoryTypes -contains $storyType

    if (-not ($isHeaderStory -or $isFooterStory)) {
        continue
    }

    $shapes = $story.InlineShapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)

        if ($isHeaderStory) {
            # BTec_Logo-Orange inline picture: image1.jpg -> image2.jpg
            $shape.Name = "image2.jpg"
        } else {
            # PearsonLogo inline picture: image2.png -> image1.png
            #
            # Footer-hosted inline shapes aren't addressable directly
            # through Range.InlineShapes for a rename (the footer story
            # has several paragraphs ahead of the picture's own
            # paragraph); selecting the shape first and renaming it via
            # the resulting Selection.InlineShapes collection reaches
            # the same picture through a path that resolves correctly.
            [void]$shape.Select()
            $word.Selection.InlineShapes.Item(1).Name = "image1.png"
        }
    }
}
